# Update the quicksort_median_of_three results table with new benchmark values.
# The BEST_CASE / AVERAGE_CASE rows swap meaning (row 2 now holds the data that
# used to be labelled AVERAGE_CASE's position etc.) and every data row receives
# freshly measured numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quicksort_median_of_three")

# Row labels: row 2 is now AVERAGE_CASE, row 3 is now BEST_CASE, row 4 stays WORST_CASE.
$ws.Range("A2").Value = "AVERAGE_CASE"
$ws.Range("A3").Value = "BEST_CASE"
$ws.Range("A4").Value = "WORST_CASE"

$row2 = @(7470, 13820, 28740, 62330, 149580, 398580, 1185890, 3916180, 14205750)
$row3 = @(1960, 3880, 7830, 16440, 34330, 71750, 151010, 313910, 657430)
$row4 = @(7530, 27650, 105990, 416240, 1653020, 6588860, 26384260, 105780140, 421945370)

for ($i = 0; $i -lt $row2.Length; $i++) {
    $col = $i + 2   # columns B..J
    $ws.Cells.Item(2, $col).Value = $row2[$i]
    $ws.Cells.Item(3, $col).Value = $row3[$i]
    $ws.Cells.Item(4, $col).Value = $row4[$i]
}
